# Update the "Price" (column D) and "Volume(1h)" (column E) figures
# on the cryptos list sheet, per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'30.484.58"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.54%  '
$ws.Cells.Item(3, 4).Value = "'2.093.91"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -1.17%  '
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = "'330.41"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.18%  '
$ws.Cells.Item(6, 5).Value = '  +0.09%  '
$ws.Cells.Item(7, 5).Value = '  -2.24%  '
$ws.Cells.Item(8, 4).Value = "'0.4428"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.05%  '
$ws.Cells.Item(9, 4).Value = "'54.03"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +15.98%  '
$ws.Cells.Item(10, 4).Value = "'0.08937"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.32%  '
$ws.Cells.Item(11, 4).Value = "'1.151"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.64%  '
$ws.Cells.Item(12, 4).Value = "'24.27"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -4.14%  '
$ws.Cells.Item(13, 4).Value = "'2.092.54"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.37%  '
$ws.Cells.Item(14, 4).Value = "'6.686"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -1.79%  '
$ws.Cells.Item(15, 4).Value = "'7.706"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -1.76%  '
$ws.Cells.Item(16, 4).Value = "'95.92"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -2.27%  '
$ws.Cells.Item(17, 5).Value = '  +0.15%  '
$ws.Cells.Item(18, 4).Value = "'0.00001121"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.75%  '
$ws.Cells.Item(19, 4).Value = "'0.06613"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.93%  '
$ws.Cells.Item(20, 4).Value = "'19.12"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.69%  '
$ws.Cells.Item(21, 5).Value = '  +0.06%  '
$ws.Cells.Item(22, 4).Value = "'6.262"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -2.51%  '
$ws.Cells.Item(23, 4).Value = "'30.509.07"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.75%  '
$ws.Cells.Item(24, 4).Value = "'12.31"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.73%  '
$ws.Cells.Item(25, 4).Value = "'2.310"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.75%  '
$ws.Cells.Item(26, 4).Value = "'2.336.80"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.36%  '
$ws.Cells.Item(27, 4).Value = "'22.23"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -3.31%  '
$ws.Cells.Item(28, 4).Value = "'2.568"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.39%  '
$ws.Cells.Item(29, 4).Value = "'163.52"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.19%  '
$ws.Cells.Item(30, 4).Value = "'131.61"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.94%  '
$ws.Cells.Item(31, 4).Value = "'1.187"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.45%  '
$ws.Cells.Item(32, 5).Value = '  -0.89%  '
$ws.Cells.Item(33, 4).Value = "'1.654"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +8.04%  '
$ws.Cells.Item(34, 4).Value = "'6.151"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -2.00%  '
$ws.Cells.Item(35, 4).Value = "'3.904"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -2.59%  '
$ws.Cells.Item(36, 4).Value = "'10.10"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.17%  '
$ws.Cells.Item(37, 4).Value = "'0.02559"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -2.82%  '
$ws.Cells.Item(38, 4).Value = "'0.06792"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.14%  '
$ws.Cells.Item(39, 4).Value = "'5.462"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.02%  '
$ws.Cells.Item(40, 4).Value = "'12.63"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -4.73%  '
$ws.Cells.Item(41, 4).Value = "'0.2253"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.05%  '
$ws.Cells.Item(42, 4).Value = "'0.6871"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.09%  '
$ws.Cells.Item(43, 5).Value = '  -0.81%  '
$ws.Cells.Item(44, 4).Value = "'1.001"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.11%  '
$ws.Cells.Item(45, 4).Value = "'13.94"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -1.31%  '
$ws.Cells.Item(46, 4).Value = "'0.6321"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.48%  '
$ws.Cells.Item(47, 4).Value = "'2.197"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.78%  '
$ws.Cells.Item(48, 4).Value = "'3.626"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.61%  '
$ws.Cells.Item(49, 5).Value = '  +5.17%  '
$ws.Cells.Item(50, 4).Value = "'1.243"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'81.56"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -2.15%  '
